# Re-pulled dSF (column F) data update.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    3  = -2
    13 = -3
    15 = -6
    18 = -3
    19 = 0
    20 = -3
    21 = -1
    24 = 0
    27 = -3
    28 = 0
    30 = -1
    37 = -3
    38 = -3
    41 = 2
    48 = -2
    56 = 0
    60 = 0
}

foreach ($row in $updates.Keys) {
    $ws.Range("F$row").Value = $updates[$row]
}
